$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$computers = @("Computer1","Computer2","Computer3","Computer4","Computer5","Computer6","Computer7","Computer8","Computer9")
$tests = @("Test1","Test2","Test3","Test4","Test5","Test6","Test7","Test8","Test9")
$quantities = @(10,20,30,40,50,60,70,80,90)

for ($i = 0; $i -lt 9; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $computers[$i]
    $ws.Cells.Item($row, 2).Value = $tests[$i]
    $ws.Cells.Item($row, 3).Value = 200
    $ws.Cells.Item($row, 4).Value = $quantities[$i]
}

$ws.Range("B11").Select()
